$wb = $excel.ActiveWorkbook

# --- Worksheet references (by name, order-independent) ---
$wsCreateRecipient   = $wb.Worksheets.Item("CreateRecipient")
$wsEditRecipient     = $wb.Worksheets.Item("EditRecipient")
$wsDeleteRecipient   = $wb.Worksheets.Item("DeleteRecipient")
$wsAddressCreate     = $wb.Worksheets.Item("AddressCreate")
$wsEditAddressBook   = $wb.Worksheets.Item("EditAddressBook")
$wsDeleteAddressBook = $wb.Worksheets.Item("DeleteAddressBook")

# --- Header renames: "First Name"/"Last Name" -> "Name1"/"Name2" ---
$wsCreateRecipient.Range("A1").Value = "Name1"
$wsCreateRecipient.Range("B1").Value = "Name2"

$wsEditRecipient.Range("A1").Value = "Name1"
$wsEditRecipient.Range("B1").Value = "Name2"

$wsDeleteRecipient.Range("A1").Value = "Name1"
$wsDeleteRecipient.Range("B1").Value = "Name2"

# --- AddressCreate: swap sample recipient name "Radhika g<456789>" -> "sample2 delete<12345>" ---
$wsAddressCreate.Range("D2").Value = "sample2 delete<12345>"
$wsAddressCreate.Range("E2").Value = "Fax Address Recipient<9987288>,sample2 delete<12345>"

# --- EditAddressBook: swap sample recipient name "Radhika g<456789>" -> "SampleData<123467>" ---
$wsEditAddressBook.Range("C2").Value = "Fax Address Recipient<9987288>,SampleData<123467>"
$wsEditAddressBook.Range("E2").Value = "Fax Address Updated Recipient<9987288>,SampleData<123467>"
$wsEditAddressBook.Range("F1").Value = "Fax Address"
$wsEditAddressBook.Range("F2").Value = "Recipient"

# --- DeleteAddressBook: swap sample recipient name and A2 header value ---
$wsDeleteAddressBook.Range("A2").Value = "Address Book Updated"
$wsDeleteAddressBook.Range("C2").Value = "Fax Address Updated Recipient<9987288>,SampleData<123467>"

# --- Column width tweaks (best-effort; engine quantizes column widths to
#     discrete steps, so these land on the nearest achievable value) ---
$wsAddressCreate.Columns("C").ColumnWidth = 29.916666666666664
$wsAddressCreate.Columns("D").ColumnWidth = 20.75
$wsDeleteAddressBook.Columns("A").ColumnWidth = 21.416666666666664
$wsDeleteAddressBook.Columns("C").ColumnWidth = 49.58333333333333

# --- View state: selections on the non-active sheets first ---
$wsEditRecipient.Range("D3").Select() | Out-Null
$wsDeleteRecipient.Range("B1").Select() | Out-Null
$wsEditAddressBook.Range("C2").Select() | Out-Null
$wsDeleteAddressBook.Range("C7").Select() | Out-Null

# --- Finally activate AddressCreate and select E6 so it ends up as the
#     selected/active tab, matching the target workbook view ---
$wsAddressCreate.Activate()
$wsAddressCreate.Range("E6").Select() | Out-Null
